$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 40; this shifts the existing rows 40-72 down
# to 41-73 (row insert copies formatting from the row above/below per
# Excel's default behaviour, matching the style="2" date format already
# present on column D for these data rows).
$ws.Rows.Item(40).Insert()

# Populate the newly inserted row 40 with the new data record.
$ws.Range("A40").Value = 11
$ws.Range("B40").Value = "Vega Monumental Concepción"
$ws.Range("C40").Value = "Bíobío"
$ws.Range("D40").Value = 44586
$ws.Range("E40").Value = 8
$ws.Range("F40").Value = 100112024
$ws.Range("G40").Value = "Choclo"
$ws.Range("H40").Value = "Choclero"
$ws.Range("I40").Value = "Primera"
$ws.Range("J40").Value = 30000
$ws.Range("K40").Value = 150
$ws.Range("L40").Value = 200
$ws.Range("M40").Value = 175
$ws.Range("N40").Value = "$/unidad"
$ws.Range("O40").Value = "Región Metropolitana"
$ws.Range("P40").Value = 175
$ws.Range("Q40").Value = 1
$ws.Range("R40").Value = "Hortaliza"
